# Reorders the comma-separated "Recorded By" values in column G.
# Rule observed from the target diff:
#   - If the literal token "System" (exact case) is present among the
#     comma-separated entries, move it to the end of the list (other
#     entries keep their relative order).
#   - Otherwise, rotate the list left by one position (the first entry
#     moves to the end).
#   - Single-value cells (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasSystem = $true }
    }

    $newParts = @()
    if ($hasSystem) {
        $systemCount = 0
        foreach ($p in $parts) {
            if ($p.Equals("System")) {
                $systemCount = $systemCount + 1
            } else {
                $newParts += $p
            }
        }
        for ($i = 0; $i -lt $systemCount; $i++) {
            $newParts += "System"
        }
    } else {
        for ($i = 1; $i -lt $parts.Length; $i++) {
            $newParts += $parts[$i]
        }
        $newParts += $parts[0]
    }

    $newVal = [string]::Join(", ", $newParts)
    if (-not $newVal.Equals($val)) {
        $cell.Value2 = $newVal
    }
}
